# Correct the filtering for product list and Transaction list
#
# The data table (A1:F33, header in row 1) gets sorted ascending by
# column B ("name"), keeping each row's data (and formatting) together,
# and the active selection moves from E34 to D5 (scrolled back to the
# top of the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sort the data range (excluding the header row) by column B, ascending.
$dataRange = $ws.Range("A2:F33")
$sortKey   = $ws.Range("B2:B33")
$dataRange.Sort($sortKey, 1)

# Move / collapse the selection to D5 (also resets any scrolled
# "topLeftCell" back to the default top-left of the sheet).
$ws.Range("D5").Select()
